# B6-PowerPoint.pptx edit
#  1. Re-style the three summary tables (slides 14-16) from the old
#     custom table style to the built-in "No Style, No Grid" style.
#  2. Re-point the deck's theme color scheme from the "Integral" /
#     "Red Violet" palette to the stock "Office" palette.

$p = $ppt.ActivePresentation

function HexToComRgb($hex) {
    # PowerPoint's RGB property is a COLORREF (0x00BBGGRR), i.e. the
    # byte order is reversed relative to the usual "RRGGBB" hex string.
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

# --- 1. Tables: switch style id on every table shape on slides 14-16 ---
$newTableStyleId = "{55A4179E-852A-4A9E-BBCD-0810D7ADCCDB}"
$tableSlideIndexes = @(14, 15, 16)
foreach ($slideIdx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Theme: swap the "Integral" color scheme for the "Office" one ---
$officeColors = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToComRgb $officeColors[$i]
}
